# "Brainstorm websites to scrape"
#
# The document ends with a paragraph that holds only the (empty)
# "_GoBack" bookmark. This script:
#   1. inserts a new blank paragraph right before that bookmark paragraph
#   2. types "5. Websites for weather data" into the bookmark paragraph
#      (ahead of the bookmark tags, so the bookmark still brackets the
#      very end of the paragraph)
#   3. appends two more paragraphs after it, describing (a) and (b)
#      candidate weather-data websites.

$d = $word.ActiveDocument

# The last paragraph in the document is the one holding bookmarkStart/
# bookmarkEnd for "_GoBack".
$lastIndex = $d.Paragraphs.Count
$goBackPara = $d.Paragraphs.Item($lastIndex)

# 1. Blank paragraph inserted immediately before the bookmark paragraph.
$goBackPara.Range.InsertParagraphBefore()

# Re-fetch: the bookmark paragraph shifted down by one.
$lastIndex = $d.Paragraphs.Count
$goBackPara = $d.Paragraphs.Item($lastIndex)

# 2. Put the heading text into the (still-empty) bookmark paragraph,
# inserting it *before* the collapsed range so it lands ahead of the
# bookmark start/end markers rather than after them.
$insertPoint = $d.Range($goBackPara.Range.Start, $goBackPara.Range.Start)
$insertPoint.InsertBefore("weather data")

$insertPoint2 = $d.Range($goBackPara.Range.Start, $goBackPara.Range.Start)
$insertPoint2.InsertBefore("5. Websites for ")

# 3. Two new paragraphs after the heading paragraph.
$lastIndex = $d.Paragraphs.Count
$goBackPara = $d.Paragraphs.Item($lastIndex)
$goBackPara.Range.InsertParagraphAfter()

$aIndex = $lastIndex + 1
$aPara = $d.Paragraphs.Item($aIndex)
$aInsertPoint = $d.Range($aPara.Range.Start, $aPara.Range.Start)
$aInsertPoint.InsertBefore("(a) Use Wunderground; Yahoo weather API ")

$aPara = $d.Paragraphs.Item($aIndex)
$aPara.Range.InsertParagraphAfter()

$bIndex = $aIndex + 1
$bPara = $d.Paragraphs.Item($bIndex)
$bInsertPoint = $d.Range($bPara.Range.Start, $bPara.Range.Start)
$bInsertPoint.InsertBefore("(b) Realtime science grade feeds: NOAA radar composites, NASA MODIS (weird format to scrape since all are projected images, but worth a look) ")
